# "hopefully final commit, integrated matrix ops"
#
# The 2x3 matrix in B1:C3 is replaced by a single 1x4 row vector in B1:E1:
#   B1=2, C1=2 (unchanged), plus two new values D1=4.4, E1=1231.2
# Rows 2 and 3 (B2:C3) are removed entirely.
# The sheet's default column width and the current selection are also updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 2 and 3 completely - the matrix now lives entirely on row 1.
$ws.Rows("2:3").Delete()

# Extend row 1 with the new matrix entries (B1:C1 already hold 2 and 2).
$ws.Range("D1").Value = 4.4
$ws.Range("E1").Value = 1231.2

# Sheet-wide default column width grows slightly.
$ws.StandardWidth = 11.640625

# Move the active selection/cursor to I7.
$ws.Range("I7").Select()
